# Update column G ("K" - strikeouts) values for rows 2-44.
# The save_data regeneration switched the source stat from "Strike#" to "K",
# so the K column values below are recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,1,1,1,3,2,0,0,2,2,2,0,4,1,0,1,2,2,2,1,0,5,2,4,1,2,1,4,0,0,1,0,4,1,1,1,2,3,1,2,0,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
